# Tripadvisor New Orleans shard 143 - restructure workbook:
#   - Sheet in tab position 1 (currently "hotel_info") becomes "review_info":
#       header row gets the full 25-column review_info header, data row removed.
#   - Sheet in tab position 2 (currently "review_info") becomes "hotel_info":
#       header row truncated to 10 columns with a new "State" column inserted,
#       and a data row for the Marriott New Orleans hotel (incl. State) is added.

$wb = $excel.ActiveWorkbook
$sheetAtPos1 = $wb.Worksheets.Item(1)
$sheetAtPos2 = $wb.Worksheets.Item(2)

# --- Rename sheets (swap names), avoiding name collisions ---
$sheetAtPos2.Name = "hotel_info_tmp_swap"
$sheetAtPos1.Name = "review_info"
$sheetAtPos2.Name = "hotel_info"

$reviewSheet = $sheetAtPos1   # tab position 1, now named "review_info"
$hotelSheet  = $sheetAtPos2   # tab position 2, now named "hotel_info"

# --- Rebuild review_info sheet: full header row, no data rows ---
$reviewHeaders = @(
    "STR","reviewer_ID","reviewer_name","Review_ID","Date_of_scraping","ReviewURL",
    "Tripadvisor_gcode","Tripadvisor_dcode","Tripadvisor_rcode","review_date","review_title",
    "review_content","review_rating","trip_month","trip_purpose","value","rooms","Location",
    "Cleanliness","Sleep Quality","Service","Picture(yes=1)","respondent","response_date","response_text"
)

for ($i = 0; $i -lt $reviewHeaders.Length; $i++) {
    $reviewSheet.Cells.Item(1, $i + 1).Value = $reviewHeaders[$i]
}

# Remove the old hotel data row (row 2) that used to live on this sheet.
$reviewSheet.Rows.Item(2).Delete()

# --- Rebuild hotel_info sheet: 10-column header (State inserted) + 1 data row ---
$hotelHeaders = @("STR","Hotel_Name","State","City","Zip","TA_ReviewURL","Tripadvisor_Hotel_Name","English_Reviews_num","Local_Rank","Total_Reviews_num")

for ($i = 0; $i -lt $hotelHeaders.Length; $i++) {
    $hotelSheet.Cells.Item(1, $i + 1).Value = $hotelHeaders[$i]
}

# Drop the now-unused trailing columns (K1:Y1) left over from the old review_info header.
$hotelSheet.Range("K1:Y1").Clear()

# Data row (row 2): hotel record, with new "State" value inserted.
$hotelSheet.Cells.Item(2, 1).Value = 6821
$hotelSheet.Cells.Item(2, 2).Value = "Marriott New Orleans"
$hotelSheet.Cells.Item(2, 3).Value = "Louisiana"
$hotelSheet.Cells.Item(2, 4).Value = "New Orleans"
$hotelSheet.Cells.Item(2, 5).Value = 70130
$hotelSheet.Cells.Item(2, 6).Value = "https://www.tripadvisor.com/Hotel_Review-g60864-d89101-Reviews-New_Orleans_Marriott-New_Orleans_Louisiana.html"
$hotelSheet.Cells.Item(2, 7).Value = "New Orleans Marriott"

# These look numeric but must stay text, like in the source data - force text storage.
$hotelSheet.Cells.Item(2, 8).NumberFormat = "@"
$hotelSheet.Cells.Item(2, 8).Value = "3555"
$hotelSheet.Cells.Item(2, 9).NumberFormat = "@"
$hotelSheet.Cells.Item(2, 9).Value = "88"
$hotelSheet.Cells.Item(2, 10).NumberFormat = "@"
$hotelSheet.Cells.Item(2, 10).Value = "3684"
